$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1304.1666
$ws.Range("I28").Value = 1208.3704
$ws.Range("K28").Value = 1208.3704
$ws.Range("M28").Value = -723.3704
$ws.Range("H70").Value = 3909.2222
$ws.Range("I70").Value = 3400.6667
$ws.Range("J70").Value = 4163.5
$ws.Range("K70").Value = 10202.0001
$ws.Range("L70").Value = 12490.5
$ws.Range("M70").Value = -9932.000100000001
$ws.Range("N70").Value = -13030.5
$ws.Range("H73").Value = 3909.2222
$ws.Range("I73").Value = 3400.6667
$ws.Range("J73").Value = 4163.5
$ws.Range("K73").Value = 10202.0001
$ws.Range("L73").Value = 12490.5
$ws.Range("M73").Value = -9266.000100000001
$ws.Range("N73").Value = -14362.5
$ws.Range("H92").Value = 4808958
$ws.Range("I92").Value = 1017.4211
$ws.Range("K92").Value = 1017.4211
$ws.Range("M92").Value = 230.5789
$ws.Range("H96").Value = 250.83333
$ws.Range("I96").Value = 228.66667
$ws.Range("K96").Value = 686.00001
$ws.Range("M96").Value = 686.99999
$ws.Range("H127").Value = 9602.132
$ws.Range("I127").Value = 1718.4286
$ws.Range("J127").Value = 19340.824
$ws.Range("K127").Value = 5155.2858
$ws.Range("L127").Value = 58022.472
$ws.Range("M127").Value = -195.2857999999997
$ws.Range("N127").Value = -67942.47200000001
$ws.Range("H132").Value = 679.54236
$ws.Range("I132").Value = 675.98114
$ws.Range("K132").Value = 2027.94342
$ws.Range("M132").Value = 502.0565799999999
$ws.Range("H138").Value = 4429.5835
$ws.Range("J138").Value = 4556.9565
$ws.Range("L138").Value = 13670.8695
$ws.Range("N138").Value = -23950.8695

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12013.483
$ws.Range("I32").Value = 8964.468999999999
$ws.Range("J32").Value = 36100.7
$ws.Range("K32").Value = 8964.468999999999
$ws.Range("L32").Value = 36100.7
$ws.Range("M32").Value = -8677.468999999999
$ws.Range("N32").Value = -36674.7
$ws.Range("H63").Value = 7899
$ws.Range("J63").Value = 7899
$ws.Range("L63").Value = 7899
$ws.Range("N63").Value = -9271
$ws.Range("H66").Value = 7899
$ws.Range("J66").Value = 7899
$ws.Range("L66").Value = 39495
$ws.Range("N66").Value = -46359
$ws.Range("H74").Value = 2138.037
$ws.Range("I74").Value = 1830.5416
$ws.Range("J74").Value = 4598
$ws.Range("K74").Value = 1830.5416
$ws.Range("L74").Value = 4598
$ws.Range("M74").Value = -956.5416
$ws.Range("N74").Value = -6346
$ws.Range("H77").Value = 2138.037
$ws.Range("I77").Value = 1830.5416
$ws.Range("J77").Value = 4598
$ws.Range("K77").Value = 9152.708000000001
$ws.Range("L77").Value = 22990
$ws.Range("M77").Value = -4784.708000000001
$ws.Range("N77").Value = -31726
$ws.Range("H102").Value = 33334574
$ws.Range("I102").Value = 1377.6666
$ws.Range("J102").Value = 333333340
$ws.Range("K102").Value = 1377.6666
$ws.Range("L102").Value = 333333340
$ws.Range("M102").Value = 244.3334
$ws.Range("N102").Value = -333336584
$ws.Range("H110").Value = 3125.077
$ws.Range("I110").Value = 3119.652
$ws.Range("J110").Value = 3166.6667
$ws.Range("K110").Value = 3119.652
$ws.Range("L110").Value = 3166.6667
$ws.Range("M110").Value = -1074.652
$ws.Range("N110").Value = -7256.6667
$ws.Range("H122").Value = 3898.16
$ws.Range("I122").Value = 2864.1667
$ws.Range("J122").Value = 6557
$ws.Range("K122").Value = 8592.500100000001
$ws.Range("L122").Value = 19671
$ws.Range("M122").Value = -6142.500100000001
$ws.Range("N122").Value = -24571
$ws.Range("H132").Value = 3619.468
$ws.Range("I132").Value = 3320.7954
$ws.Range("J132").Value = 8000
$ws.Range("K132").Value = 9962.386200000001
$ws.Range("L132").Value = 24000
$ws.Range("M132").Value = -7432.386200000001
$ws.Range("N132").Value = -29060

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 12501780
$ws.Range("I94").Value = 4763033.5
$ws.Range("K94").Value = 4763033.5
$ws.Range("M94").Value = -4762582.5
$ws.Range("H105").Value = 2875.1191
$ws.Range("J105").Value = 5015.273
$ws.Range("L105").Value = 5015.273
$ws.Range("N105").Value = -8509.273000000001
$ws.Range("H107").Value = 777.4706
$ws.Range("I107").Value = 765.5714
$ws.Range("K107").Value = 765.5714
$ws.Range("M107").Value = 1154.4286

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 551
$ws.Range("I16").Value = 469.66666
$ws.Range("K16").Value = 469.66666
$ws.Range("M16").Value = -182.66666
$ws.Range("H31").Value = 4020.9424
$ws.Range("I31").Value = 3666.2856
$ws.Range("J31").Value = 4151.6055
$ws.Range("K31").Value = 3666.2856
$ws.Range("L31").Value = 4151.6055
$ws.Range("M31").Value = -3371.2856
$ws.Range("N31").Value = -4741.6055
$ws.Range("H34").Value = 4020.9424
$ws.Range("I34").Value = 3666.2856
$ws.Range("J34").Value = 4151.6055
$ws.Range("K34").Value = 3666.2856
$ws.Range("L34").Value = 4151.6055
$ws.Range("M34").Value = -3464.2856
$ws.Range("N34").Value = -4555.6055
$ws.Range("H58").Value = 3136.4285
$ws.Range("I58").Value = 1318.3334
$ws.Range("K58").Value = 1318.3334
$ws.Range("M58").Value = -1115.3334
$ws.Range("H62").Value = 5559.143
$ws.Range("J62").Value = 2750
$ws.Range("L62").Value = 2750
$ws.Range("N62").Value = -3998
$ws.Range("H65").Value = 5559.143
$ws.Range("J65").Value = 2750
$ws.Range("L65").Value = 13750
$ws.Range("N65").Value = -19990
$ws.Range("H99").Value = 18705938
$ws.Range("I99").Value = 4884950
$ws.Range("J99").Value = 28578072
$ws.Range("K99").Value = 4884950
$ws.Range("L99").Value = 28578072
$ws.Range("M99").Value = -4883452
$ws.Range("N99").Value = -28581068
$ws.Range("H105").Value = 488.5
$ws.Range("I105").Value = 504.5
$ws.Range("J105").Value = 464.5
$ws.Range("K105").Value = 504.5
$ws.Range("L105").Value = 464.5
$ws.Range("M105").Value = 1242.5
$ws.Range("N105").Value = -3958.5
$ws.Range("H113").Value = 551
$ws.Range("I113").Value = 469.66666
$ws.Range("K113").Value = 469.66666
$ws.Range("M113").Value = 1700.33334
$ws.Range("H122").Value = 332543.03
$ws.Range("I122").Value = 681735
$ws.Range("K122").Value = 2045205
$ws.Range("M122").Value = -2042755
$ws.Range("H126").Value = 18705938
$ws.Range("I126").Value = 4884950
$ws.Range("J126").Value = 28578072
$ws.Range("K126").Value = 14654850
$ws.Range("L126").Value = 85734216
$ws.Range("M126").Value = -14652380
$ws.Range("N126").Value = -85739156
$ws.Range("H132").Value = 3481.0527
$ws.Range("I132").Value = 3510.8462
$ws.Range("J132").Value = 3416.5
$ws.Range("K132").Value = 10532.5386
$ws.Range("L132").Value = 10249.5
$ws.Range("M132").Value = -8002.5386
$ws.Range("N132").Value = -15309.5
$ws.Range("H134").Value = 3454.3462
$ws.Range("I134").Value = 2455.65
$ws.Range("J134").Value = 6783.3335
$ws.Range("K134").Value = 7366.950000000001
$ws.Range("L134").Value = 20350.0005
$ws.Range("M134").Value = -4831.950000000001
$ws.Range("N134").Value = -25420.0005
$ws.Range("H136").Value = 3136.4285
$ws.Range("I136").Value = 1318.3334
$ws.Range("K136").Value = 3955.0002
$ws.Range("M136").Value = -1405.0002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 7253.524
$ws.Range("I56").Value = 7253.524
$ws.Range("K56").Value = 7253.524
$ws.Range("M56").Value = -6723.524
$ws.Range("H122").Value = 1519.619
$ws.Range("J122").Value = 1702
$ws.Range("L122").Value = 15318
$ws.Range("N122").Value = -20218
$ws.Range("H132").Value = 3265.1
$ws.Range("J132").Value = 5259.6
$ws.Range("L132").Value = 47336.4
$ws.Range("N132").Value = -52396.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 17858062
$ws.Range("I97").Value = 950.8889
$ws.Range("J97").Value = 50000864
$ws.Range("K97").Value = 950.8889
$ws.Range("L97").Value = 50000864
$ws.Range("M97").Value = -454.8889
$ws.Range("N97").Value = -50001856
$ws.Range("H113").Value = 5486.6
$ws.Range("I113").Value = 4415.2856
$ws.Range("K113").Value = 4415.2856
$ws.Range("M113").Value = -2245.2856
$ws.Range("H122").Value = 5470.8125
$ws.Range("I122").Value = 2964.875
$ws.Range("J122").Value = 7976.75
$ws.Range("K122").Value = 8894.625
$ws.Range("L122").Value = 23930.25
$ws.Range("M122").Value = -6444.625
$ws.Range("N122").Value = -28830.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2792.6128
$ws.Range("I61").Value = 1772.7826
$ws.Range("K61").Value = 1772.7826
$ws.Range("M61").Value = -1570.7826
$ws.Range("H113").Value = 2792.6128
$ws.Range("I113").Value = 1772.7826
$ws.Range("K113").Value = 1772.7826
$ws.Range("M113").Value = 397.2174

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 12348535
$ws.Range("J62").Value = 18521560
$ws.Range("L62").Value = 18521560
$ws.Range("N62").Value = -18522808
$ws.Range("H65").Value = 12348535
$ws.Range("J65").Value = 18521560
$ws.Range("L65").Value = 92607800
$ws.Range("N65").Value = -92614040
$ws.Range("H107").Value = 468.06668
$ws.Range("I107").Value = 401.6154
$ws.Range("K107").Value = 1204.8462
$ws.Range("M107").Value = 715.1538
$ws.Range("H113").Value = 553.3333
$ws.Range("I113").Value = 581.4211
$ws.Range("K113").Value = 1744.2633
$ws.Range("M113").Value = 425.7366999999999
$ws.Range("H122").Value = 3906.087
$ws.Range("I122").Value = 3290
$ws.Range("K122").Value = 9870
$ws.Range("M122").Value = -7420
$ws.Range("H126").Value = 2046
$ws.Range("I126").Value = 1624
$ws.Range("K126").Value = 4872
$ws.Range("M126").Value = -2402
